$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells in this sheet are stored as text (prices/percentages
# with leading zeros / trailing zeros that must be preserved verbatim), so
# force a text number format before writing each value to avoid Excel
# auto-converting numeric-looking strings (e.g. "1.00" -> 1, "6.50" -> 6.5).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.183.23'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.302.84'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.18%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.54'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.68'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.34%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.302.20'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.21%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.72%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.869.14'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.16%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.95'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.296.07'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.81%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.185.61'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.21'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.41'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.64'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '374.23'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.79%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.542'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.45%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.444.77'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.29%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -7.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.172'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.89%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.56'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.74%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '165.62'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.23%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.73'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.39%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.52'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.44%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'RenzoRestakedETH'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.331.65'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.33%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.73'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -13.51%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.99'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.752'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.19'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.05%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.03%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.368.80'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.88%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.50'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.25'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.20%  '
